$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6848073333333332
$ws.Range("H2").Value = 2.054422
$ws.Range("I2").Value = 0.2268310526442471
$ws.Range("J2").Value = 0.2268310526442472
$ws.Range("M2").Value = 3.558321333333333
$ws.Range("N2").Value = 10.674964
$ws.Range("O2").Value = 0.3039644761000113
$ws.Range("P2").Value = 0.3039644761000113
$ws.Range("Q2").Value = 2.436764543423111
$ws.Range("R2").Value = 21.93088089080799
$ws.Range("S2").Value = 0.06894858208022266
$ws.Range("T2").Value = 0.06894858208022266
$ws.Range("G3").Value = 0.6848073333333332
$ws.Range("H3").Value = 2.054422
$ws.Range("I3").Value = 0.2268310526442471
$ws.Range("J3").Value = 0.2268310526442472
$ws.Range("M3").Value = 5.383140666666667
$ws.Range("O3").Value = 0.4598470400038817
$ws.Range("P3").Value = 0.4598470400038817
$ws.Range("Q3").Value = 3.686414204898222
$ws.Range("R3").Value = 33.177727844084
$ws.Range("S3").Value = 0.1043075881394217
$ws.Range("T3").Value = 0.1043075881394217
$ws.Range("G4").Value = 0.6848073333333332
$ws.Range("H4").Value = 2.054422
$ws.Range("I4").Value = 0.2268310526442471
$ws.Range("J4").Value = 0.2268310526442472
$ws.Range("M4").Value = 2.764910333333333
$ws.Range("N4").Value = 8.294730999999999
$ws.Range("O4").Value = 0.2361884838961071
$ws.Range("P4").Value = 0.236188483896107
$ws.Range("Q4").Value = 1.893430872275777
$ws.Range("R4").Value = 17.04087785048199
$ws.Range("S4").Value = 0.05357488242460277
$ws.Range("T4").Value = 0.05357488242460277
$ws.Range("I5").Value = 0.1086184939966157
$ws.Range("J5").Value = 0.1086184939966157
$ws.Range("M5").Value = 3.558321333333333
$ws.Range("N5").Value = 10.674964
$ws.Range("O5").Value = 0.3039644761000113
$ws.Range("P5").Value = 0.3039644761000113
$ws.Range("Q5").Value = 1.166849476055111
$ws.Range("R5").Value = 10.501645284496
$ws.Range("S5").Value = 0.03301616362245351
$ws.Range("T5").Value = 0.0330161636224535
$ws.Range("I6").Value = 0.1086184939966157
$ws.Range("J6").Value = 0.1086184939966157
$ws.Range("M6").Value = 5.383140666666667
$ws.Range("O6").Value = 0.4598470400038817
$ws.Range("P6").Value = 0.4598470400038817
$ws.Range("S6").Value = 0.04994789295402311
$ws.Range("T6").Value = 0.04994789295402311
$ws.Range("I7").Value = 0.1086184939966157
$ws.Range("J7").Value = 0.1086184939966157
$ws.Range("M7").Value = 2.764910333333333
$ws.Range("N7").Value = 8.294730999999999
$ws.Range("O7").Value = 0.2361884838961071
$ws.Range("P7").Value = 0.236188483896107
$ws.Range("Q7").Value = 0.9066730830537777
$ws.Range("R7").Value = 8.160057747483998
$ws.Range("S7").Value = 0.02565443742013906
$ws.Range("T7").Value = 0.02565443742013906
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.006290666666667
$ws.Range("H8").Value = 6.018872
$ws.Range("I8").Value = 0.6645504533591371
$ws.Range("J8").Value = 0.6645504533591372
$ws.Range("M8").Value = 3.558321333333333
$ws.Range("N8").Value = 10.674964
$ws.Range("O8").Value = 0.3039644761000113
$ws.Range("P8").Value = 0.3039644761000113
$ws.Range("Q8").Value = 7.139026880067555
$ws.Range("R8").Value = 64.251241920608
$ws.Range("S8").Value = 0.2019997303973351
$ws.Range("T8").Value = 0.2019997303973351
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.006290666666667
$ws.Range("H9").Value = 6.018872
$ws.Range("I9").Value = 0.6645504533591371
$ws.Range("J9").Value = 0.6645504533591372
$ws.Range("M9").Value = 5.383140666666667
$ws.Range("O9").Value = 0.4598470400038817
$ws.Range("P9").Value = 0.4598470400038817
$ws.Range("Q9").Value = 10.80014487688711
$ws.Range("R9").Value = 97.20130389198401
$ws.Range("S9").Value = 0.3055915589104368
$ws.Range("T9").Value = 0.3055915589104369
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.006290666666667
$ws.Range("H10").Value = 6.018872
$ws.Range("I10").Value = 0.6645504533591371
$ws.Range("J10").Value = 0.6645504533591372
$ws.Range("M10").Value = 2.764910333333333
$ws.Range("N10").Value = 8.294730999999999
$ws.Range("O10").Value = 0.2361884838961071
$ws.Range("P10").Value = 0.236188483896107
$ws.Range("Q10").Value = 5.547213795936888
$ws.Range("R10").Value = 49.924924163432
$ws.Range("S10").Value = 0.1569591640513652
$ws.Range("T10").Value = 0.1569591640513652
